# Apply scheduled-runner profit recalculation updates to the Leve profit
# columns (currentAveragePrice / LevePrice* / LeveProfit*) across sheets.
$wb = $excel.ActiveWorkbook

# ALC row 32
$ws = $wb.Sheets.Item("ALC")
$ws.Range("H32").Value = 17996.234
$ws.Range("J32").Value = 14003.637
$ws.Range("L32").Value = 14003.637
$ws.Range("N32").Value = -14655.637

# ALC row 38
$ws = $wb.Sheets.Item("ALC")
$ws.Range("H38").Value = 282.1
$ws.Range("I38").Value = 227.11111
$ws.Range("J38").Value = 777
$ws.Range("K38").Value = 681.3333299999999
$ws.Range("L38").Value = 2331
$ws.Range("M38").Value = -309.3333299999999
$ws.Range("N38").Value = -3075

# ALC row 58
$ws = $wb.Sheets.Item("ALC")
$ws.Range("H58").Value = 2123.111
$ws.Range("I58").Value = 390
$ws.Range("J58").Value = 3509.6
$ws.Range("K58").Value = 1170
$ws.Range("L58").Value = 10528.8
$ws.Range("M58").Value = -1020
$ws.Range("N58").Value = -10828.8

# ALC row 125
$ws = $wb.Sheets.Item("ALC")
$ws.Range("H125").Value = 10000.667
$ws.Range("J125").Value = 10000
$ws.Range("L125").Value = 90000
$ws.Range("N125").Value = -94920

# ALC row 132
$ws = $wb.Sheets.Item("ALC")
$ws.Range("H132").Value = 323738.34
$ws.Range("I132").Value = 363790.88
$ws.Range("K132").Value = 1091372.64
$ws.Range("M132").Value = -1088842.64

# ALC row 135
$ws = $wb.Sheets.Item("ALC")
$ws.Range("H135").Value = 6514.091
$ws.Range("I135").Value = 3753.0908
$ws.Range("J135").Value = 12036.091
$ws.Range("K135").Value = 33777.8172
$ws.Range("L135").Value = 108324.819
$ws.Range("M135").Value = -31242.8172
$ws.Range("N135").Value = -113394.819

# ALC row 137
$ws = $wb.Sheets.Item("ALC")
$ws.Range("H137").Value = 6043
$ws.Range("I137").Value = 4285.885
$ws.Range("K137").Value = 12857.655
$ws.Range("M137").Value = -10307.655

# ARM row 63
$ws = $wb.Sheets.Item("ARM")
$ws.Range("H63").Value = 5466.1665
$ws.Range("J63").Value = 6799.5
$ws.Range("L63").Value = 6799.5
$ws.Range("N63").Value = -8171.5

# ARM row 66
$ws = $wb.Sheets.Item("ARM")
$ws.Range("H66").Value = 5466.1665
$ws.Range("J66").Value = 6799.5
$ws.Range("L66").Value = 33997.5
$ws.Range("N66").Value = -40861.5

# ARM row 132
$ws = $wb.Sheets.Item("ARM")
$ws.Range("H132").Value = 545912.9399999999
$ws.Range("J132").Value = 117781.25
$ws.Range("L132").Value = 353343.75
$ws.Range("N132").Value = -358403.75

# BSM row 134
$ws = $wb.Sheets.Item("BSM")
$ws.Range("H134").Value = 1931365.4
$ws.Range("I134").Value = 2639687.2
$ws.Range("K134").Value = 7919061.600000001
$ws.Range("M134").Value = -7916526.600000001

# CRP row 31
$ws = $wb.Sheets.Item("CRP")
$ws.Range("H31").Value = 4883.625
$ws.Range("I31").Value = 1198.4166
$ws.Range("J31").Value = 8568.833000000001
$ws.Range("K31").Value = 1198.4166
$ws.Range("L31").Value = 8568.833000000001
$ws.Range("M31").Value = -903.4166
$ws.Range("N31").Value = -9158.833000000001

# CRP row 34
$ws = $wb.Sheets.Item("CRP")
$ws.Range("H34").Value = 4883.625
$ws.Range("I34").Value = 1198.4166
$ws.Range("J34").Value = 8568.833000000001
$ws.Range("K34").Value = 1198.4166
$ws.Range("L34").Value = 8568.833000000001
$ws.Range("M34").Value = -996.4166
$ws.Range("N34").Value = -8972.833000000001

# CRP row 58
$ws = $wb.Sheets.Item("CRP")
$ws.Range("H58").Value = 23267416
$ws.Range("I58").Value = 40006892
$ws.Range("J58").Value = 18143.889
$ws.Range("K58").Value = 40006892
$ws.Range("L58").Value = 18143.889
$ws.Range("M58").Value = -40006689
$ws.Range("N58").Value = -18549.889

# CRP row 62
$ws = $wb.Sheets.Item("CRP")
$ws.Range("H62").Value = 7028.4443
$ws.Range("I62").Value = 8375.333000000001
$ws.Range("J62").Value = 4334.6665
$ws.Range("K62").Value = 8375.333000000001
$ws.Range("L62").Value = 4334.6665
$ws.Range("M62").Value = -7751.333000000001
$ws.Range("N62").Value = -5582.6665

# CRP row 65
$ws = $wb.Sheets.Item("CRP")
$ws.Range("H65").Value = 7028.4443
$ws.Range("I65").Value = 8375.333000000001
$ws.Range("J65").Value = 4334.6665
$ws.Range("K65").Value = 41876.665
$ws.Range("L65").Value = 21673.3325
$ws.Range("M65").Value = -38756.665
$ws.Range("N65").Value = -27913.3325

# CRP row 74
$ws = $wb.Sheets.Item("CRP")
$ws.Range("H74").Value = 24899.8
$ws.Range("J74").Value = 24899.8
$ws.Range("L74").Value = 24899.8
$ws.Range("N74").Value = -26647.8

# CRP row 77
$ws = $wb.Sheets.Item("CRP")
$ws.Range("H77").Value = 24899.8
$ws.Range("J77").Value = 24899.8
$ws.Range("L77").Value = 74699.39999999999
$ws.Range("N77").Value = -83435.39999999999

# CRP row 122
$ws = $wb.Sheets.Item("CRP")
$ws.Range("H122").Value = 2789.6
$ws.Range("I122").Value = 2217.8
$ws.Range("K122").Value = 6653.400000000001
$ws.Range("M122").Value = -4203.400000000001

# CRP row 134
$ws = $wb.Sheets.Item("CRP")
$ws.Range("H134").Value = 62506390
$ws.Range("I134").Value = 83338240
$ws.Range("K134").Value = 250014720
$ws.Range("M134").Value = -250012185

# CRP row 136
$ws = $wb.Sheets.Item("CRP")
$ws.Range("H136").Value = 23267416
$ws.Range("I136").Value = 40006892
$ws.Range("J136").Value = 18143.889
$ws.Range("K136").Value = 120020676
$ws.Range("L136").Value = 54431.667
$ws.Range("M136").Value = -120018126
$ws.Range("N136").Value = -59531.667

# CUL row 37
$ws = $wb.Sheets.Item("CUL")
$ws.Range("H37").Value = 99563.21000000001
$ws.Range("J37").Value = 99563.21000000001
$ws.Range("L37").Value = 298689.63
$ws.Range("N37").Value = -298913.63

# CUL row 132
$ws = $wb.Sheets.Item("CUL")
$ws.Range("H132").Value = 3330
$ws.Range("J132").Value = 3996.4285
$ws.Range("L132").Value = 35967.8565
$ws.Range("N132").Value = -41027.8565

# GSM row 126
$ws = $wb.Sheets.Item("GSM")
$ws.Range("H126").Value = 16137009
$ws.Range("I126").Value = 26319968
$ws.Range("K126").Value = 78959904
$ws.Range("M126").Value = -78957434

# LTW row 16
$ws = $wb.Sheets.Item("LTW")
$ws.Range("H16").Value = 2170.0417
$ws.Range("I16").Value = 755.0454999999999
$ws.Range("K16").Value = 755.0454999999999
$ws.Range("M16").Value = -585.0454999999999

# LTW row 61
$ws = $wb.Sheets.Item("LTW")
$ws.Range("H61").Value = 5851.364
$ws.Range("I61").Value = 4757.4346
$ws.Range("J61").Value = 8367.4
$ws.Range("K61").Value = 4757.4346
$ws.Range("L61").Value = 8367.4
$ws.Range("M61").Value = -4555.4346
$ws.Range("N61").Value = -8771.4

# LTW row 68
$ws = $wb.Sheets.Item("LTW")
$ws.Range("H68").Value = 3072.5454
$ws.Range("I68").Value = 2271.4285
$ws.Range("J68").Value = 4474.5
$ws.Range("K68").Value = 2271.4285
$ws.Range("L68").Value = 4474.5
$ws.Range("M68").Value = -1522.4285
$ws.Range("N68").Value = -5972.5

# LTW row 71
$ws = $wb.Sheets.Item("LTW")
$ws.Range("H71").Value = 3072.5454
$ws.Range("I71").Value = 2271.4285
$ws.Range("J71").Value = 4474.5
$ws.Range("K71").Value = 11357.1425
$ws.Range("L71").Value = 22372.5
$ws.Range("M71").Value = -7613.1425
$ws.Range("N71").Value = -29860.5

# LTW row 113
$ws = $wb.Sheets.Item("LTW")
$ws.Range("H113").Value = 5851.364
$ws.Range("I113").Value = 4757.4346
$ws.Range("J113").Value = 8367.4
$ws.Range("K113").Value = 4757.4346
$ws.Range("L113").Value = 8367.4
$ws.Range("M113").Value = -2587.4346
$ws.Range("N113").Value = -12707.4

# WVR row 132
$ws = $wb.Sheets.Item("WVR")
$ws.Range("H132").Value = 9557.842000000001
$ws.Range("I132").Value = 9300
$ws.Range("K132").Value = 27900
$ws.Range("M132").Value = -25370

# WVR row 138
$ws = $wb.Sheets.Item("WVR")
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
